$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from 2023-09-14 (45183) to 2023-09-15 (45184)
$ws.Range("C2:C7").Value = 45184
